$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "quantidade"

# Update codes (was food names, now letter codes)
$ws.Range("A2").Value = "A"
$ws.Range("A3").Value = "B"
$ws.Range("A4").Value = "C"

# Move selection to reflect final cursor position when saved
$ws.Range("J8").Select()
